# The author re-typed the "word" text for two rows, wrapping the letter
# in quotation marks. Excel's AutoCorrect turns straight quotes into
# curly/smart quotes (U+201C "left double quotation mark" and
# U+201D "right double quotation mark") as you type, so the stored
# strings end up with those Unicode characters rather than a plain ".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$leftQuote = [char]0x201C
$rightQuote = [char]0x201D

# Row 29 (B29) was edited first ("Words with a P" -> Words with a "P"),
# then row 5 (B5) ("Consonant K" -> Consonant "K"). Keeping this order
# matches where the two new shared-string entries land at the end of
# the shared-strings table.
$ws.Range("B29").Value = "Words with a " + $leftQuote + "P" + $rightQuote
$ws.Range("B5").Value = "Consonant " + $leftQuote + "K" + $rightQuote

# Cursor ends up on B6 after the edits.
$ws.Range("B6").Select()
